# [rdbms] - remove double-commit scenario to avoid unnecessary error condition.
# - allow for non-standard or non-CRUD statements.
#
# This reproduces the edits made to the hidden "#system" lookup sheet that
# backs the workbook's named ranges of Nexial command names:
#
#   - "base"     (column F) loses the "clearVariables(variables)" entry at
#                 F19 (it duplicated "clear(variables)"), so every entry
#                 below it shifts up one row and the range shrinks from
#                 F2:F40 to F2:F39.
#   - "external" (column J) gains a new "terminate(programName)" entry at
#                 J6, growing the range from J2:J5 to J2:J6.
#   - "web"      (column Z) gains two new entries, "saveSelectedText(var,locator)"
#                 and "saveSelectedValue(var,locator)", inserted before the
#                 old Z99 ("saveTableAsCsv..."), so everything from the old
#                 Z99 downward shifts down two rows and the range grows from
#                 Z2:Z135 to Z2:Z137.
#
# NOTE: Range.Delete()/Range.Insert() on this host shift the *entire row*
# (every column), not just the target column, which would corrupt the other
# independent lookup columns (H, L, etc.) that happen to share those rows.
# So the column shifts below are done by copying single-column cell values
# directly, which only touches column F (or column Z) and leaves every
# other column on those rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- "base" (column F): remove the "clearVariables(variables)" row (F19). ---
# Shift F20:F40 up into F19:F39, then blank out the now-unused F40.
for ($r = 19; $r -le 39; $r++) {
    $srcAddr = "F" + ($r + 1)
    $dstAddr = "F" + $r
    $ws.Range($dstAddr).Value = $ws.Range($srcAddr).Value()
}
$ws.Range("F40").ClearContents()

# --- "external" (column J): append "terminate(programName)" at J6. ---
$ws.Range("J6").Value = "terminate(programName)"

# --- "web" (column Z): insert two new rows before the old Z99 and fill them. ---
# Shift Z99:Z135 down into Z101:Z137 (walk bottom-up so writes don't clobber
# values still waiting to be copied), then fill the freed Z99:Z100.
for ($r = 135; $r -ge 99; $r--) {
    $srcAddr = "Z" + $r
    $dstAddr = "Z" + ($r + 2)
    $ws.Range($dstAddr).Value = $ws.Range($srcAddr).Value()
}
$ws.Range("Z99").Value = "saveSelectedText(var,locator)"
$ws.Range("Z100").Value = "saveSelectedValue(var,locator)"

# --- Update the defined names (they are not auto-resized by the above). ---
$wb.Names.Item("base").RefersTo = "='#system'!`$F`$2:`$F`$39"
$wb.Names.Item("external").RefersTo = "='#system'!`$J`$2:`$J`$6"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$137"
